$d = $word.ActiveDocument

# The "Learning Agreement" table (3rd table in the document) lists, per
# semester, the host-institution course taken and the home course it is
# exempted against, along with credit counts. Update the course codes
# and credit numbers for semesters 1 and 2, and clear out the row that
# had erroneously been filled in for semester 3.

$t = $d.Tables.Item(3)

# Semester 1 row
$t.Cell(3, 2).Range.Text = "1"
$t.Cell(3, 3).Range.Text = "EE571"
$t.Cell(3, 4).Range.Text = "1"
$t.Cell(3, 5).Range.Text = "CS315"
$t.Cell(3, 6).Range.Text = "3"

# Semester 2 row
$t.Cell(4, 2).Range.Text = "2"
$t.Cell(4, 3).Range.Text = "CS124"
$t.Cell(4, 4).Range.Text = "2"
$t.Cell(4, 5).Range.Text = "CS224"
$t.Cell(4, 6).Range.Text = "4"

# Semester 3 row — bugfix: clear the values that were mistakenly entered
$t.Cell(5, 2).Range.Text = ""
$t.Cell(5, 3).Range.Text = ""
$t.Cell(5, 4).Range.Text = ""
$t.Cell(5, 5).Range.Text = ""
$t.Cell(5, 6).Range.Text = ""
